# Weekly data refresh: two new price records (week of 2021-11-05, date
# serial 44505) are inserted right after row 424, pushing the remaining
# historical rows down by two positions (old A425:R491 -> A427:R493).
#
# The two new rows are seeded as copies of the (soon to be former) row
# 425/426 content and then have their Fecha/Volumen/Precio fields
# overwritten with the new week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 425/426; everything below shifts down by two.
$ws.Rows.Item(425).Insert()
$ws.Rows.Item(426).Insert()

# Seed the new rows with the data currently one above-shift's-worth below
# them (i.e. what used to be rows 425/426, now sitting at 427/428) so all
# the columns that don't change (Mercado, Región, Calidad, Unidad, Origen,
# etc.) come along for free.
$ws.Range("A427:R427").Copy($ws.Range("A425"))
$ws.Range("A428:R428").Copy($ws.Range("A426"))

# Row 425 (Primera): new Fecha + updated price figures.
$ws.Cells.Item(425, 4).Value = 44505
$ws.Cells.Item(425, 10).Value = 16100
$ws.Cells.Item(425, 12).Value = 600
$ws.Cells.Item(425, 13).Value = 552
$ws.Cells.Item(425, 16).Value = 552

# Row 426 (Segunda): new Fecha + updated price figures.
$ws.Cells.Item(426, 4).Value = 44505
$ws.Cells.Item(426, 10).Value = 4900
$ws.Cells.Item(426, 12).Value = 400
$ws.Cells.Item(426, 13).Value = 400
$ws.Cells.Item(426, 16).Value = 400
